# Auto-generated edit script: updates market-price derived columns (H:N)
# across multiple worksheets to match the scheduled runner's refreshed values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2764.1428
$ws.Range("I18").Value = 2489.8
$ws.Range("J18").Value = 3450
$ws.Range("K18").Value = 2489.8
$ws.Range("L18").Value = 3450
$ws.Range("M18").Value = -2205.8
$ws.Range("N18").Value = -4018
$ws.Range("H33").Value = 154.6923
$ws.Range("I33").Value = 91.90909000000001
$ws.Range("K33").Value = 91.90909000000001
$ws.Range("M33").Value = 137.09091
$ws.Range("H34").Value = 4098.875
$ws.Range("I34").Value = 4098.875
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 4098.875
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3895.875
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 4098.875
$ws.Range("I36").Value = 4098.875
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4098.875
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -3383.875
$ws.Range("N36").ClearContents()
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H53").Value = 1763.4286
$ws.Range("I53").Value = 93.666664
$ws.Range("J53").Value = 3015.75
$ws.Range("K53").Value = 93.666664
$ws.Range("L53").Value = 3015.75
$ws.Range("M53").Value = 543.333336
$ws.Range("N53").Value = -4289.75
$ws.Range("H64").Value = 2800
$ws.Range("J64").Value = 3200
$ws.Range("L64").Value = 3200
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 2800
$ws.Range("J67").Value = 3200
$ws.Range("L67").Value = 3200
$ws.Range("N67").Value = -4916
$ws.Range("H98").Value = 1715.9474
$ws.Range("I98").Value = 1981.6428
$ws.Range("K98").Value = 1981.6428
$ws.Range("M98").Value = -483.6428000000001
$ws.Range("H103").Value = 700
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H112").Value = 1860
$ws.Range("I112").Value = 1050
$ws.Range("J112").Value = 2130
$ws.Range("K112").Value = 3150
$ws.Range("L112").Value = 6390
$ws.Range("M112").Value = -2042
$ws.Range("N112").Value = -8606
$ws.Range("H122").Value = 1715.9474
$ws.Range("I122").Value = 1981.6428
$ws.Range("K122").Value = 5944.928400000001
$ws.Range("M122").Value = -3494.928400000001
$ws.Range("H125").Value = 975.2308
$ws.Range("I125").Value = 972.5454999999999
$ws.Range("K125").Value = 8752.9095
$ws.Range("M125").Value = -6292.9095
$ws.Range("H132").Value = 2292.5715
$ws.Range("I132").Value = 2389.6
$ws.Range("J132").Value = 2050
$ws.Range("K132").Value = 7168.799999999999
$ws.Range("L132").Value = 6150
$ws.Range("M132").Value = -4638.799999999999
$ws.Range("N132").Value = -11210
$ws.Range("H138").Value = 5118.256
$ws.Range("I138").Value = 2866.6875
$ws.Range("J138").Value = 11668.272
$ws.Range("K138").Value = 8600.0625
$ws.Range("L138").Value = 35004.81600000001
$ws.Range("M138").Value = -3460.0625
$ws.Range("N138").Value = -45284.81600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2567.2307
$ws.Range("I26").Value = 910.5714
$ws.Range("K26").Value = 910.5714
$ws.Range("M26").Value = -580.5714
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 3000
$ws.Range("K61").Value = 3000
$ws.Range("M61").Value = -2788
$ws.Range("H74").Value = 1092.5834
$ws.Range("I74").Value = 1055.5454
$ws.Range("K74").Value = 1055.5454
$ws.Range("M74").Value = -181.5454
$ws.Range("H77").Value = 1092.5834
$ws.Range("I77").Value = 1055.5454
$ws.Range("K77").Value = 5277.727
$ws.Range("M77").Value = -909.7269999999999
$ws.Range("H95").Value = 31912.666
$ws.Range("J95").Value = 31912.666
$ws.Range("L95").Value = 31912.666
$ws.Range("N95").Value = -37404.666
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774
$ws.Range("H122").Value = 3732.4443
$ws.Range("I122").Value = 4173.75
$ws.Range("J122").Value = 3379.4
$ws.Range("K122").Value = 12521.25
$ws.Range("L122").Value = 10138.2
$ws.Range("M122").Value = -10071.25
$ws.Range("N122").Value = -15038.2
$ws.Range("H132").Value = 1182.7826
$ws.Range("I132").Value = 1209.3889
$ws.Range("K132").Value = 3628.1667
$ws.Range("M132").Value = -1098.1667
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 35500
$ws.Range("J92").Value = 35500
$ws.Range("L92").Value = 35500
$ws.Range("N92").Value = -40492
$ws.Range("H95").Value = 30502.555
$ws.Range("J95").Value = 30502.555
$ws.Range("L95").Value = 30502.555
$ws.Range("N95").Value = -35994.555
$ws.Range("H99").Value = 13230.826
$ws.Range("J99").Value = 14334.667
$ws.Range("L99").Value = 14334.667
$ws.Range("N99").Value = -17330.667
$ws.Range("H126").Value = 13230.826
$ws.Range("J126").Value = 14334.667
$ws.Range("L126").Value = 43004.001
$ws.Range("N126").Value = -47944.001
$ws.Range("H132").Value = 2187.25
$ws.Range("I132").Value = 2187.25
$ws.Range("K132").Value = 6561.75
$ws.Range("M132").Value = -4031.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 611
$ws.Range("I18").Value = 533.2
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 1599.6
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -1430.6
$ws.Range("N18").Value = -3338
$ws.Range("H97").Value = 8570.666999999999
$ws.Range("J97").Value = 9154.666999999999
$ws.Range("L97").Value = 27464.001
$ws.Range("N97").Value = -28456.001
$ws.Range("H113").Value = 1181.8334
$ws.Range("I113").Value = 998
$ws.Range("J113").Value = 1218.6
$ws.Range("K113").Value = 2994
$ws.Range("L113").Value = 3655.8
$ws.Range("M113").Value = -824
$ws.Range("N113").Value = -7995.799999999999
$ws.Range("H131").Value = 1479
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1479
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 4437
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -14517

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1253
$ws.Range("I2").Value = 380
$ws.Range("J2").Value = 2999
$ws.Range("K2").Value = 380
$ws.Range("L2").Value = 2999
$ws.Range("M2").Value = -267
$ws.Range("N2").Value = -3225
$ws.Range("H11").Value = 6573571
$ws.Range("I11").Value = 6669999.5
$ws.Range("J11").Value = 6501249.5
$ws.Range("K11").Value = 6669999.5
$ws.Range("L11").Value = 6501249.5
$ws.Range("M11").Value = -6669860.5
$ws.Range("N11").Value = -6501527.5
$ws.Range("H98").Value = 14127
$ws.Range("J98").Value = 14127
$ws.Range("L98").Value = 14127
$ws.Range("N98").Value = -20117

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H100").Value = 1780.6
$ws.Range("I100").Value = 1725.75
$ws.Range("K100").Value = 1725.75
$ws.Range("M100").Value = -1184.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7889
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 65299.668
$ws.Range("J64").Value = 65299.668
$ws.Range("L64").Value = 65299.668
$ws.Range("N64").Value = -65795.66800000001
$ws.Range("H65").Value = 7889
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 65299.668
$ws.Range("J67").Value = 65299.668
$ws.Range("L67").Value = 65299.668
$ws.Range("N67").Value = -67015.66800000001
$ws.Range("H94").Value = 46994.6
$ws.Range("I94").Value = 49993.5
$ws.Range("J94").Value = 34999
$ws.Range("K94").Value = 49993.5
$ws.Range("L94").Value = 34999
$ws.Range("M94").Value = -49092.5
$ws.Range("N94").Value = -36801
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

Write-Host "Applied market price updates across sheets."